$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.290242
$ws.Range("H2").Value = 0.870726
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 33.35034833333334
$ws.Range("N2").Value = 100.051045
$ws.Range("O2").Value = 0.2287589433580892
$ws.Range("P2").Value = 0.2287589433580892
$ws.Range("Q2").Value = 9.679671800963334
$ws.Range("R2").Value = 87.11704620867
$ws.Range("S2").Value = 0.2287589433580892
$ws.Range("T2").Value = 0.2287589433580892

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.290242
$ws.Range("H3").Value = 0.870726
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 25.677964
$ws.Range("N3").Value = 77.03389199999999
$ws.Range("O3").Value = 0.176132010781908
$ws.Range("P3").Value = 0.176132010781908
$ws.Range("Q3").Value = 7.452823627288
$ws.Range("R3").Value = 67.075412645592
$ws.Range("S3").Value = 0.176132010781908
$ws.Range("T3").Value = 0.176132010781908

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.290242
$ws.Range("H4").Value = 0.870726
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 73.55277366666667
$ws.Range("N4").Value = 220.658321
$ws.Range("O4").Value = 0.5045181122808869
$ws.Range("P4").Value = 0.5045181122808869
$ws.Range("Q4").Value = 21.34810413456067
$ws.Range("R4").Value = 192.132937211046
$ws.Range("S4").Value = 0.5045181122808869
$ws.Range("T4").Value = 0.5045181122808869

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.290242
$ws.Range("H5").Value = 0.870726
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.20708666666667
$ws.Range("N5").Value = 39.62126
$ws.Range("O5").Value = 0.09059093357911581
$ws.Range("P5").Value = 0.09059093357911581
$ws.Range("Q5").Value = 3.833251248306667
$ws.Range("R5").Value = 34.49926123476
$ws.Range("S5").Value = 0.09059093357911581
$ws.Range("T5").Value = 0.09059093357911581
